$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final 6-row table: 2 sending clusters (FAPs, sCs) x 3 target clusters (ECs, FAPs, sCs)
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#          E..T numeric metrics

$rows = @(
    @("FAPs","Gdf7","Bmpr1a","ECs", 2, 0.6666666666666666, 0.210033, 0.630099, 0.311901169594117, 0.311901169594117, 3, 1, 4.750436666666666, 14.25131, 0.07745299862590357, 0.07745299862590359, 0.99774846441, 8.979736179689999, 0.02415768085999086, 0.02415768085999086),
    @("FAPs","Gdf7","Bmpr1a","FAPs", 2, 0.6666666666666666, 0.210033, 0.630099, 0.311901169594117, 0.311901169594117, 3, 1, 34.05277366666667, 102.158321, 0.5552098927072401, 0.5552098927072401, 7.152206211531, 64.369855903779, 0.1731706149056124, 0.1731706149056124),
    @("FAPs","Gdf7","Bmpr1a","sCs", 2, 0.6666666666666666, 0.210033, 0.630099, 0.311901169594117, 0.311901169594117, 3, 1, 22.52994333333334, 67.58983000000001, 0.3673371086668564, 0.3673371086668564, 4.732031588130001, 42.58828429317, 0.1145728738285137, 0.1145728738285137),
    @("sCs","Gdf7","Bmpr1a","ECs", 3, 1, 0.463363, 1.390089, 0.6880988304058829, 0.6880988304058829, 3, 1, 4.750436666666666, 14.25131, 0.07745299862590357, 0.07745299862590359, 2.201176585176666, 19.81058926659, 0.05329531776591271, 0.05329531776591272),
    @("sCs","Gdf7","Bmpr1a","FAPs", 3, 1, 0.463363, 1.390089, 0.6880988304058829, 0.6880988304058829, 3, 1, 34.05277366666667, 102.158321, 0.5552098927072401, 0.5552098927072401, 15.77879536450767, 142.009158280569, 0.3820392778016277, 0.3820392778016277),
    @("sCs","Gdf7","Bmpr1a","sCs", 3, 1, 0.463363, 1.390089, 0.6880988304058829, 0.6880988304058829, 3, 1, 22.52994333333334, 67.58983000000001, 0.3673371086668564, 0.3673371086668564, 10.43954213276333, 93.95587919487001, 0.2527642348383426, 0.2527642348383426)
)

$rowIndex = 2
foreach ($r in $rows) {
    $colIndex = 1
    foreach ($val in $r) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }
    $rowIndex = $rowIndex + 1
}
